$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Fix title text spacing: "Analyzing Projects" -> "Analyzing Projects " (i #10)
$titleShape = $s.Shapes.Item(5)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Characters(1, 18).Text = "Analyzing Projects "

# Reposition shapes that were shifted down (vertical offset only; x/width/height unchanged)
$s.Shapes.Item(14).Top = 424.4143307086614  # id=373 "Line"
$s.Shapes.Item(47).Top = 434.45  # id=454 "Manipulate Variables"
$s.Shapes.Item(48).Top = 461.2848053496068  # id=455 "TextBox 454"
$s.Shapes.Item(49).Top = 590.7180480960573  # id=456 "Graphic 455"
$s.Shapes.Item(50).Top = 624.9290551181102  # id=457 "TextBox 456"
$s.Shapes.Item(51).Top = 550.3540344480251  # id=460 "Graphic 459"
$s.Shapes.Item(52).Top = 561.1526794653536  # id=465 "TextBox 464"
$s.Shapes.Item(53).Top = 550.1066141732283  # id=466 "Graphic 465"
$s.Shapes.Item(54).Top = 561.4402362204725  # id=470 "TextBox 469"
$s.Shapes.Item(55).Top = 555.5409448818898  # id=472 "Table"
$s.Shapes.Item(56).Top = 586.7829284858216  # id=474 "Graphic 473"
$s.Shapes.Item(57).Top = 598.5699212598425  # id=475 "TextBox 474"
$s.Shapes.Item(58).Top = 518.7921447842455  # id=476 "TextBox 475"
$s.Shapes.Item(59).Top = 611.1277952755905  # id=477 "Straight Arrow Connector 476"
$s.Shapes.Item(60).Top = 593.4886780173172  # id=478 "Straight Arrow Connector 477"
$s.Shapes.Item(61).Top = 570.8667908535367  # id=479 "Straight Arrow Connector 478"
$s.Shapes.Item(62).Top = 611.1277952755905  # id=483 "Straight Arrow Connector 482"
$s.Shapes.Item(63).Top = 571.0576477952721  # id=484 "Straight Arrow Connector 483"
$s.Shapes.Item(64).Top = 603.7686614173228  # id=486 "Straight Arrow Connector 485"
$s.Shapes.Item(65).Top = 623.0199280598401  # id=488 "Straight Arrow Connector 487"
$s.Shapes.Item(66).Top = 712.6322937645634  # id=489 "TextBox 488"
$s.Shapes.Item(68).Top = 258.9114228228359  # id=30 "Graphic 29"
$s.Shapes.Item(69).Top = 240.15047464094476  # id=32 "Graphic 31"
$s.Shapes.Item(70).Top = 278.1669291338583  # id=35 "TextBox 34"
$s.Shapes.Item(71).Top = 265.5662995125985  # id=36 "Graphic 35"
$s.Shapes.Item(72).Top = 304.3843307086614  # id=37 "TextBox 36"
$s.Shapes.Item(73).Top = 271.1270905141741  # id=38 "Table"
$s.Shapes.Item(74).Top = 286.64378362755997  # id=39 "Straight Arrow Connector 38"
$s.Shapes.Item(75).Top = 305.25338582677165  # id=41 "TextBox 40"
$s.Shapes.Item(76).Top = 231.24850393700788  # id=44 "Rectangle 43"
$s.Shapes.Item(77).Top = 244.66268166535406  # id=52 "Straight Arrow Connector 51"
$s.Shapes.Item(78).Top = 262.30181102362207  # id=53 "Straight Arrow Connector 52"
$s.Shapes.Item(79).Top = 237.09299472598408  # id=54 "Straight Arrow Connector 53"
$s.Shapes.Item(80).Top = 222.96543307086614  # id=492 "Rectangle 491"
$s.Shapes.Item(81).Top = 274.0012665425212  # id=493 "Straight Arrow Connector 492"
$s.Shapes.Item(82).Top = 246.6403937007874  # id=495 "Graphic 494"
$s.Shapes.Item(83).Top = 273.67079170157575  # id=496 "Straight Arrow Connector 495"
$s.Shapes.Item(84).Top = 352.7201690803176  # id=6 "Google Shape;66;p1"
$s.Shapes.Item(87).Top = 259.3078766157486  # id=12 "TextBox 11"
$s.Shapes.Item(88).Top = 559.3982849165295  # id=13 "TextBox 12"
$s.Shapes.Item(89).Top = 525.7120667440878  # id=15 "Graphic 14"
$s.Shapes.Item(90).Top = 613.4784251968504  # id=21 "TextBox 20"
$s.Shapes.Item(91).Top = 577.2688293976341  # id=22 "Graphic 21"
$s.Shapes.Item(92).Top = 668.1829133858267  # id=23 "TextBox 22"
$s.Shapes.Item(93).Top = 631.0352755905512  # id=24 "Graphic 23"
$s.Shapes.Item(94).Top = 550.3586120771553  # id=26 "Graphic 25"
$s.Shapes.Item(95).Top = 583.9457703314883  # id=47 "TextBox 46"
$s.Shapes.Item(96).Top = 583.3876647952663  # id=48 "TextBox 47"
$s.Shapes.Item(97).Top = 584.1076377952755  # id=49 "TextBox 48"
$s.Shapes.Item(98).Top = 570.9411316621946  # id=458 "Straight Arrow Connector 457"
$s.Shapes.Item(99).Top = 536.119293238582  # id=459 "Rectangle 458"
$s.Shapes.Item(100).Top = 633.4117322834645  # id=481 "Graphic 480"
$s.Shapes.Item(101).Top = 644.2104188007788  # id=490 "TextBox 489"
$s.Shapes.Item(102).Top = 633.1643372086592  # id=491 "Graphic 490"
$s.Shapes.Item(103).Top = 644.4979553559047  # id=494 "TextBox 493"
$s.Shapes.Item(104).Top = 638.5986634173222  # id=497 "Table"
$s.Shapes.Item(105).Top = 653.9244881889764  # id=498 "Straight Arrow Connector 497"
$s.Shapes.Item(106).Top = 654.1153543307087  # id=499 "Straight Arrow Connector 498"
$s.Shapes.Item(107).Top = 633.4162992125985  # id=500 "Graphic 499"
$s.Shapes.Item(108).Top = 667.0034645669291  # id=501 "TextBox 500"
$s.Shapes.Item(109).Top = 666.4453543307087  # id=502 "TextBox 501"
$s.Shapes.Item(110).Top = 667.1653748307016  # id=503 "TextBox 502"
$s.Shapes.Item(111).Top = 653.9988188976378  # id=504 "Straight Arrow Connector 503"
$s.Shapes.Item(112).Top = 619.1770324740073  # id=505 "Rectangle 504"
